# Generate Report for Archive
#
# The nightly localization-status report moved the two test documents from
# "Ready for handoff" to "In Translation". Update every cell that shows that
# status (the "zh-cn"/"de-de" status columns on the Overview sheet, and the
# "Status" column on each per-locale sheet), then tighten up the status
# columns' widths to fit the new, shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the per-locale status ---
foreach ($row in 2..3) {
    foreach ($col in @("E", "F")) {
        $cell = $overview.Range("$col$row")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Per-locale sheets: column C is "Status" ---
foreach ($row in 2..3) {
    $cell = $zhcn.Range("C$row")
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

foreach ($row in 2..3) {
    $cell = $dede.Range("C$row")
    if ($cell.Value2 -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- Re-fit the status columns now that the text is shorter ---
# Target character width (matches the report-generation tool's column sizing
# for the new status text); ColumnWidth is rounded to the nearest pixel by
# the host, so feed it the closest value that lands on the desired width.
$newColumnWidth = 12.576851254417766

$overview.Range("E1").ColumnWidth = $newColumnWidth
$overview.Range("F1").ColumnWidth = $newColumnWidth
$zhcn.Range("C1").ColumnWidth = $newColumnWidth
$dede.Range("C1").ColumnWidth = $newColumnWidth

Write-Host "Updated status text and column widths"
